$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.1660323612164193, 0.468872571632237,  0.3770675980855352, 0.6140583018619121, 0.6119356173368786, 15),
    @(3, 0.2261183143199825, 0.4923620343352079, 0.4821095245989872, 0.6943410722397079, 0.6812725906144649, 14),
    @(4, 0.2360243239862733, 0.5248832413428339, 0.37481815022521,   0.6122239379714011, 0.5879653258807318, 13),
    @(5, 0.3051580102940423, 0.3983834828835194, 0.2461501466550942, 0.4961352100537657, 0.4085834536877719, 12),
    @(6, 0.3404744896530312, 0.5484269867923621, 0.5429546666577961, 0.736854576329547,  0.6853721378992834, 11),
    @(7, 0.3329151699208477, 0.5578246582442683, 0.4732217363555192, 0.6879111398687473, 0.6345507824806217, 10),
    @(8, 0.2814885994813455, 0.5772087615654313, 0.5238078509972534, 0.7237457087936711, 0.7072082591282391, 9),
    @(9, 0.3213393071819032, 0.6123030108606557, 0.5290292607895575, 0.7273439769390804, 0.6975633594976256, 8),
    @(10, 0.4424481932715923, 0.4424481932715923, 0.2455568012204569, 0.4955368817963572, 0.2410306973665816, 7),
    @(11, 0.3495035403546135, 0.4324022405235847, 0.3358454641318753, 0.5795217546666176, 0.5063904494495919, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
